$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.020.88"
$ws.Range("E2").Value = "  +2.45%  "
$ws.Range("D3").Value = "3.196.38"
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.11"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.85"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.35"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("E10").Value = "  +1.14%  "
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("D12").Value = "3.746.94"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.138"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.11%  "
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "60.045.84"
$ws.Range("E16").Value = "  +2.40%  "
$ws.Range("D17").Value = "3.196.62"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.27"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.29"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.21"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "370.54"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.58%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("E23").Value = "  -1.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.66"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("E25").Value = "  +1.78%  "
$ws.Range("E26").Value = "  +4.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +1.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.48"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.13"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("E33").Value = "  +4.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.19"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.98"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.24%  "
$ws.Range("E36").Value = "  +2.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.54"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +6.01%  "
$ws.Range("D38").Value = "2.794.73"
$ws.Range("E38").Value = "  +5.04%  "
$ws.Range("E39").Value = "  +9.02%  "
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("E43").Value = "  +2.03%  "
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("E45").Value = "  +1.28%  "
$ws.Range("D46").Value = "3.236.97"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.983"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.66"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.796"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.77%  "
$ws.Range("E51").Value = "  -0.01%  "

Write-Host "Done updating cryptos list"
